# Update taxon category name capitalization in column A
# "Dragonflies & Damselflies" -> "Dragonflies & damselflies"
# "Freshwater Crabs"          -> "Freshwater crabs"
# "Freshwater Fishes"         -> "Freshwater fishes"
# "Sharks (incl. Rays & Chimaeras)" -> "Sharks (incl. rays & chimaeras)"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A11").Value = "Dragonflies & damselflies"
$ws.Range("A12").Value = "Dragonflies & damselflies"

$ws.Range("A13").Value = "Freshwater crabs"
$ws.Range("A14").Value = "Freshwater crabs"

$ws.Range("A15").Value = "Freshwater fishes"
$ws.Range("A16").Value = "Freshwater fishes"

$ws.Range("A26").Value = "Sharks (incl. rays & chimaeras)"
$ws.Range("A27").Value = "Sharks (incl. rays & chimaeras)"
